$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (column D) values that look like plain numbers need to be
# forced to Text format first, so Excel stores them as the literal price strings
# (matching the source data feed format) rather than auto-converting them to
# numeric values.
$ws.Range("D5,D7,D8,D9,D10,D11,D13,D14,D15,D16,D18,D19,D22,D23,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D38,D39,D40,D41,D42,D43,D45,D46,D47,D49,D50,D51").NumberFormat = "@"

$ws.Range("D5").Value = "309.61"
$ws.Range("D7").Value = "0.4646"
$ws.Range("D8").Value = "0.3701"
$ws.Range("D9").Value = "0.07356"
$ws.Range("D10").Value = "0.8758"
$ws.Range("D11").Value = "20.47"
$ws.Range("D13").Value = "5.355"
$ws.Range("D14").Value = "6.509"
$ws.Range("D15").Value = "91.67"
$ws.Range("D16").Value = "0.07045"
$ws.Range("D18").Value = "0.000008692"
$ws.Range("D19").Value = "1.001"
$ws.Range("D22").Value = "5.317"
$ws.Range("D23").Value = "10.57"
$ws.Range("D25").Value = "1.898"
$ws.Range("D26").Value = "151.50"
$ws.Range("D27").Value = "18.38"
$ws.Range("D28").Value = "2.151"
$ws.Range("D29").Value = "5.329"
$ws.Range("D30").Value = "115.75"
$ws.Range("D31").Value = "0.08901"
$ws.Range("D32").Value = "0.7551"
$ws.Range("D33").Value = "1.159"
$ws.Range("D34").Value = "2.920"
$ws.Range("D35").Value = "4.455"
$ws.Range("D38").Value = "0.01967"
$ws.Range("D39").Value = "2.438"
$ws.Range("D40").Value = "0.05245"
$ws.Range("D41").Value = "2.927"
$ws.Range("D42").Value = "0.5316"
$ws.Range("D43").Value = "7.175"
$ws.Range("D45").Value = "8.478"
$ws.Range("D46").Value = "0.4974"
$ws.Range("D47").Value = "10.35"
$ws.Range("D49").Value = "103.80"
$ws.Range("D50").Value = "1.667"
$ws.Range("D51").Value = "0.06296"

# Remaining updated cells (coin name/link swaps, prices that are not plain numbers,
# and the Volume(1h) percentage strings) can be set directly as text.
$ws.Range("D2").Value = "26.872.10"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "1.809.68"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +4.18%  "
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "1.837.16"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D21").Value = "26.872.72"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("D24").Value = "2.019.30"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  -5.79%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -4.81%  "
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("E39").Value = "  +4.50%  "
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -1.50%  "
